$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Mark the "HPC analysis" flag column (J) for a subset of compounds ---
# "x" rows get written first so that "x" becomes shared-string index 113
# and "yes" becomes shared-string index 114 (matches target ordering).
$xRows   = @(17, 22, 23)
$yesRows = @(2, 3, 16, 24)

foreach ($r in $xRows) {
    $ws.Cells.Item($r, 10).Value = "x"
}
foreach ($r in $yesRows) {
    $ws.Cells.Item($r, 10).Value = "yes"
}

# --- 2. Highlight the corresponding compound names in column A ---
# Build the new look (bold font + red fill) once on a clean cell, then
# propagate it with Copy/PasteSpecial(Formats) so only a single new
# cell style gets created (instead of one per intermediate state).
$highlightRows = @(2, 3, 16, 17, 22, 23, 24)

$seed = $ws.Cells.Item($highlightRows[0], 1)
$seed.Font.Bold = $true
$seed.Interior.Color = 255
$seed.Copy()

for ($i = 1; $i -lt $highlightRows.Count; $i++) {
    $ws.Cells.Item($highlightRows[$i], 1).PasteSpecial(-4122)
}

# --- 3. Restore the selection left behind by the editing session ---
$ws.Range("M21").Select()
